# Rename ObjTables document/table declaration attributes to lowerCamelCase
# (e.g. Type= -> type=, TableName= -> tableName=, ObjTablesVersion= -> objTablesVersion=)
# across the three worksheets that embed these declarations as plain text
# in cell A1 (and A2 on the "Table of contents" sheet).

$wb = $excel.ActiveWorkbook

# --- "!!_Table of contents" sheet ---
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableName='Table of contents' description='Table/model and column/attribute definitions' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# --- "!!Child" sheet ---
$wsChild = $wb.Worksheets.Item("!!Child")
$wsChild.Range("A1").Value = "!!ObjTables type='Data' id='Child' name='Child' description='Represents a child' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# --- "!!Parent" sheet ---
$wsParent = $wb.Worksheets.Item("!!Parent")
$wsParent.Range("A1").Value = "!!ObjTables type='Data' id='Parent' name='Parent' description='Represents a parent' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"
